$wb = $excel.ActiveWorkbook

# ALC row 10
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 26333.334
$ws.Range("I10").Value = 35000
$ws.Range("J10").Value = 9000
$ws.Range("K10").Value = 35000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = -34707
$ws.Range("N10").Value = -9586

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4932.143
$ws.Range("J17").Value = 4932.143
$ws.Range("L17").Value = 14796.429
$ws.Range("N17").Value = -15132.429

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 105.86364
$ws.Range("I33").Value = 106.375
$ws.Range("J33").Value = 104.5
$ws.Range("K33").Value = 106.375
$ws.Range("L33").Value = 104.5
$ws.Range("M33").Value = 122.625
$ws.Range("N33").Value = -562.5

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2043.2632
$ws.Range("I40").Value = 2052.2222
$ws.Range("K40").Value = 2052.2222
$ws.Range("M40").Value = -1877.2222

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4667.364
$ws.Range("I64").Value = 3542.5
$ws.Range("J64").Value = 5310.143
$ws.Range("K64").Value = 3542.5
$ws.Range("L64").Value = 5310.143
$ws.Range("M64").Value = -3294.5
$ws.Range("N64").Value = -5806.143

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4667.364
$ws.Range("I67").Value = 3542.5
$ws.Range("J67").Value = 5310.143
$ws.Range("K67").Value = 3542.5
$ws.Range("L67").Value = 5310.143
$ws.Range("M67").Value = -2684.5
$ws.Range("N67").Value = -7026.143

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 16044105
$ws.Range("J112").Value = 17046800
$ws.Range("L112").Value = 51140400
$ws.Range("N112").Value = -51142616

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1016766.25
$ws.Range("J116").Value = 2770.3333
$ws.Range("L116").Value = 2770.3333
$ws.Range("N116").Value = -9654.3333

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1013.05
$ws.Range("I129").Value = 337.5
$ws.Range("J129").Value = 1088.1111
$ws.Range("K129").Value = 1012.5
$ws.Range("L129").Value = 3264.3333
$ws.Range("M129").Value = 3987.5
$ws.Range("N129").Value = -13264.3333

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1332
$ws.Range("I137").Value = 1278.2
$ws.Range("J137").Value = 1466.5
$ws.Range("K137").Value = 3834.6
$ws.Range("L137").Value = 4399.5
$ws.Range("M137").Value = -1284.6
$ws.Range("N137").Value = -9499.5

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 9824162
$ws.Range("I138").Value = 2275424.5
$ws.Range("J138").Value = 14708639
$ws.Range("K138").Value = 6826273.5
$ws.Range("L138").Value = 44125917
$ws.Range("M138").Value = -6821133.5
$ws.Range("N138").Value = -44136197

# ARM row 16
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 4
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 283
$ws.Range("N16").ClearContents()

# ARM row 46
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 12990
$ws.Range("J46").Value = 12990
$ws.Range("L46").Value = 12990
$ws.Range("N46").Value = -13628

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9748.6
$ws.Range("I74").Value = 1087.5
$ws.Range("K74").Value = 1087.5
$ws.Range("M74").Value = -213.5

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 9748.6
$ws.Range("I77").Value = 1087.5
$ws.Range("K77").Value = 5437.5
$ws.Range("M77").Value = -1069.5

# CRP row 2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 21002.4
$ws.Range("I2").Value = 21002.4
$ws.Range("K2").Value = 21002.4
$ws.Range("M2").Value = -20889.4

# CRP row 25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# CRP row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 15052
$ws.Range("I59").Value = 10104
$ws.Range("K59").Value = 10104
$ws.Range("M59").Value = -8959

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10820.357
$ws.Range("I4").Value = 95.05882
$ws.Range("J4").Value = 27395.818
$ws.Range("K4").Value = 285.17646
$ws.Range("L4").Value = 82187.454
$ws.Range("M4").Value = -173.17646
$ws.Range("N4").Value = -82411.454

# CUL row 8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1044.7273
$ws.Range("I8").Value = 1044.7273
$ws.Range("K8").Value = 3134.1819
$ws.Range("M8").Value = -2995.1819

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2508.1765
$ws.Range("I139").Value = 2477.4375
$ws.Range("K139").Value = 7432.3125
$ws.Range("M139").Value = -2292.3125

# GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 307
$ws.Range("I5").Value = 307
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 307
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -195
$ws.Range("N5").ClearContents()

# GSM row 23
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# GSM row 31
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2677
$ws.Range("I31").Value = 2677
$ws.Range("K31").Value = 2677
$ws.Range("M31").Value = -2385

# GSM row 37
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 2677
$ws.Range("I37").Value = 2677
$ws.Range("K37").Value = 2677
$ws.Range("M37").Value = -2400

# GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 29993.334
$ws.Range("I57").Value = 29990
$ws.Range("K57").Value = 29990
$ws.Range("M57").Value = -29170

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8597.9
$ws.Range("I70").Value = 19940
$ws.Range("J70").Value = 5762.375
$ws.Range("K70").Value = 19940
$ws.Range("L70").Value = 5762.375
$ws.Range("M70").Value = -19670
$ws.Range("N70").Value = -6302.375

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8597.9
$ws.Range("I73").Value = 19940
$ws.Range("J73").Value = 5762.375
$ws.Range("K73").Value = 19940
$ws.Range("L73").Value = 5762.375
$ws.Range("M73").Value = -19004
$ws.Range("N73").Value = -7634.375

# LTW row 33
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 30000
$ws.Range("I33").Value = 30000
$ws.Range("K33").Value = 30000
$ws.Range("M33").Value = -29710

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2482.6274
$ws.Range("I132").Value = 1794.1628
$ws.Range("J132").Value = 6183.125
$ws.Range("K132").Value = 5382.4884
$ws.Range("L132").Value = 18549.375
$ws.Range("M132").Value = -2852.4884
$ws.Range("N132").Value = -23609.375
